$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits alone in the empty
#    ListParagraph right before the contributions table.
# ---------------------------------------------------------------------------
try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
    Write-Host "Step 1 (remove _GoBack): " $_.Exception.Message
}

# ---------------------------------------------------------------------------
# 2) In the "Gozal Alizada" task-description cell, insert a leading space
#    as its own run before the existing "Introduction, Program Design..."
#    text (same sz/szCs formatting).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Introduction, Program Design (DFD), Sharing link on social media.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $pos2 = $rng2.Start
    $ins2 = $d.Range($pos2, $pos2)
    $ins2.InsertBefore(" ")
    # Nudge formatting on the inserted space and restore it so the host
    # breaks the new space off into its own run instead of merging it
    # into the following run.
    $space2 = $d.Range($pos2, $pos2 + 1)
    $space2.Font.Bold = 1
    $space2.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# 3) In the "Asima Kochariyeva" task-description cell, insert a leading
#    space run followed by a fresh "_GoBack" bookmark, right before the
#    existing "Hardware and Software Specification..." text.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Hardware and Software Specification, Preparation Presentation, Data Storage Design(ERD), Sharing link on social media.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $pos3 = $rng3.Start
    $ins3 = $d.Range($pos3, $pos3)
    $ins3.InsertBefore(" ")
    $space3 = $d.Range($pos3, $pos3 + 1)
    $space3.Font.Bold = 1
    $space3.Font.Bold = 0

    $bmRange = $d.Range($pos3 + 1, $pos3 + 1)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Host "Edits applied."
